$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.308.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.839.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6257'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07370'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2890'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.834.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6631'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001045'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.239'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.319.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '236.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.74%  '
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.234'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.69%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '157.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.409'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1335'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.07120'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.471'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.477'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.016'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.022'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.151'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.783'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6885'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.579'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01822'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("E38").Value = '  -1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.233.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.730'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9439'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.996.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.926'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.682'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.850'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1128'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3876'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.48%  '
